$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (matches source
# data that is stored as text, e.g. "28.476.76", "0.00001080", "1.002"),
# then restore the default "Normal" style so no stray number format
# is left behind on the cell.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "28.501.54"
Set-TextValue "E2" "  +2.62%  "

# Row 3
Set-TextValue "D3" "1.829.55"
Set-TextValue "E3" "  +2.07%  "

# Row 4
Set-TextValue "D4" "1.002"
Set-TextValue "E4" "  -0.02%  "

# Row 5
Set-TextValue "D5" "317.33"
Set-TextValue "E5" "  +0.43%  "

# Row 6
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  -0.07%  "

# Row 7
Set-TextValue "D7" "0.5054"
Set-TextValue "E7" "  -5.47%  "

# Row 8
Set-TextValue "D8" "0.3917"
Set-TextValue "E8" "  +2.07%  "

# Row 9
Set-TextValue "D9" "0.07715"
Set-TextValue "E9" "  +3.88%  "

# Row 10
Set-TextValue "D10" "41.98"
Set-TextValue "E10" "  +1.29%  "

# Row 11
Set-TextValue "D11" "1.113"
Set-TextValue "E11" "  +2.69%  "

# Row 12
Set-TextValue "D12" "21.02"
Set-TextValue "E12" "  +3.51%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "6.253"
Set-TextValue "E13" "  +0.76%  "

# Row 14
$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D14" "1.001"
Set-TextValue "E14" "  -0.07%  "

# Row 15
Set-TextValue "D15" "7.545"
Set-TextValue "E15" "  +1.43%  "

# Row 16
Set-TextValue "D16" "1.826.28"
Set-TextValue "E16" "  +1.79%  "

# Row 17
Set-TextValue "D17" "93.49"
Set-TextValue "E17" "  +5.82%  "

# Row 18
Set-TextValue "D18" "0.00001082"
Set-TextValue "E18" "  +2.27%  "

# Row 19
Set-TextValue "D19" "0.06611"
Set-TextValue "E19" "  +1.30%  "

# Row 20
Set-TextValue "D20" "17.74"
Set-TextValue "E20" "  +2.50%  "

# Row 21
Set-TextValue "D21" "1.001"
Set-TextValue "E21" "  +0.02%  "

# Row 22
Set-TextValue "D22" "6.113"
Set-TextValue "E22" "  +2.56%  "

# Row 23
Set-TextValue "D23" "28.539.07"
Set-TextValue "E23" "  +2.61%  "

# Row 24
Set-TextValue "D24" "11.12"
Set-TextValue "E24" "  +0.16%  "

# Row 25
Set-TextValue "D25" "2.254"
Set-TextValue "E25" "  +7.57%  "

# Row 26
Set-TextValue "D26" "156.93"
Set-TextValue "E26" "  +0.12%  "

# Row 27
Set-TextValue "D27" "20.61"
Set-TextValue "E27" "  +2.00%  "

# Row 28
Set-TextValue "D28" "2.035.41"
Set-TextValue "E28" "  +1.77%  "

# Row 29
Set-TextValue "D29" "2.423"
Set-TextValue "E29" "  +4.12%  "

# Row 30
Set-TextValue "D30" "125.04"
Set-TextValue "E30" "  +3.06%  "

# Row 31
Set-TextValue "D31" "1.131"
Set-TextValue "E31" "  +1.87%  "

# Row 32
Set-TextValue "D32" "0.1088"
Set-TextValue "E32" "  -0.40%  "

# Row 33
Set-TextValue "D33" "5.656"
Set-TextValue "E33" "  +2.79%  "

# Row 34
Set-TextValue "D34" "3.665"
Set-TextValue "E34" "  +0.33%  "

# Row 35
Set-TextValue "D35" "0.07099"
Set-TextValue "E35" "  +1.78%  "

# Row 36
Set-TextValue "D36" "0.2213"
Set-TextValue "E36" "  +1.02%  "

# Row 37
Set-TextValue "D37" "9.006"
Set-TextValue "E37" "  +7.19%  "

# Row 38
Set-TextValue "D38" "0.02321"
Set-TextValue "E38" "  +2.07%  "

# Row 39
Set-TextValue "D39" "5.124"
Set-TextValue "E39" "  +1.41%  "

# Row 40
Set-TextValue "D40" "0.6232"
Set-TextValue "E40" "  +2.19%  "

# Row 41
Set-TextValue "E41" "  -1.43%  "

# Row 42
Set-TextValue "D42" "1.188"
Set-TextValue "E42" "  +2.51%  "

# Row 43
Set-TextValue "D43" "1.001"
Set-TextValue "E43" "  -0.04%  "

# Row 44
Set-TextValue "D44" "1.396"
Set-TextValue "E44" "  -1.08%  "

# Row 45
Set-TextValue "D45" "13.39"
Set-TextValue "E45" "  +0.76%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.5889"
Set-TextValue "E46" "  +3.34%  "

# Row 47
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D47" "3.714"
Set-TextValue "E47" "  +0.99%  "

# Row 48
Set-TextValue "D48" "124.23"
Set-TextValue "E48" "  -0.64%  "

# Row 49
Set-TextValue "E49" "  +3.25%  "

# Row 50
Set-TextValue "D50" "1.183"
Set-TextValue "E50" "  +1.23%  "

# Row 51
Set-TextValue "D51" "0.06927"
Set-TextValue "E51" "  +2.09%  "
